$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C53 previously held the text "NA" - the script run cleared it to a blank
# (but still present) cell, matching the other empty-string cells in column C.
# Writing a bare apostrophe forces an empty text entry rather than truly
# deleting the cell; resetting the style afterwards drops the quote-prefix
# tag so the cell stays on the default style like its column peers.
$ws.Range("C53").Value = "'"
$ws.Range("C53").Style = "Normal"

# Helper to write an ISO date-like string as literal text (not an Excel date
# serial). A leading apostrophe forces text entry; resetting the style back
# to "Normal" afterwards drops the quote-prefix style Excel would otherwise
# leave tagged on the cell, keeping it on the default style like its peers.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# New rows appended by the latest script run.
Set-TextValue $ws.Range("A54") "2025-04-22"
$ws.Range("B54").Value = "développement durable"
$ws.Range("C54").Value = 152
$ws.Range("D54").Value = 1

Set-TextValue $ws.Range("A55") "2025-04-22"
$ws.Range("B55").Value = "développement durable"
$ws.Range("C55").Value = 153
$ws.Range("D55").Value = 1

Set-TextValue $ws.Range("A56") "2025-04-22"
$ws.Range("B56").Value = "développement durable"
$ws.Range("C56").Value = 154
$ws.Range("D56").Value = 1
